$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.525.78"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "'2.324.25"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'511.76"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").Value = "'131.56"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").Value = "'0.100"
$ws.Range("E9").Value = "  -3.69%  "

$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'23.53"
$ws.Range("E13").Value = "  -1.26%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'2.738.29"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").Value = "'56.510.46"

$ws.Range("E16").Value = "  -1.89%  "

$ws.Range("D17").Value = "'2.331.82"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").Value = "'10.40"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "'327.62"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("E23").Value = "  +1.49%  "

$ws.Range("D24").Value = "'8.59"
$ws.Range("E24").Value = "  +8.65%  "

$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("D28").Value = "'167.68"
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("E29").Value = "  -3.81%  "

$ws.Range("E30").Value = "  -4.38%  "

$ws.Range("D31").Value = "'6.12"
$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("D32").Value = "'18.33"
$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("E37").Value = "  -3.93%  "

$ws.Range("D38").Value = "'38.52"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").Value = "'148.81"
$ws.Range("E40").Value = "  +7.84%  "

$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("D43").Value = "'277.26"
$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("E44").Value = "  -3.91%  "

$ws.Range("D45").Value = "'0.0929"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("E46").Value = "  -2.22%  "

$ws.Range("D47").Value = "'0.557"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").Value = "'18.15"
$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D51").Value = "'17.01"
$ws.Range("E51").Value = "  +0.93%  "
